$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns per latest crypto snapshot ---
# Price values are written as text (matching the sheet's existing text-formatted
# Price/Volume columns): set NumberFormat to Text before assignment, then restore the
# cell's normal style so no stray number formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.914.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.585.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.93%  "

$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.478"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.75%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.805.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.587.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "

$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.887.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0726"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "59.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.68%  "

$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.75%  "

$ws.Range("E22").Value = "  -1.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.15%  "

$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  -0.46%  "

$ws.Range("E29").Value = "  -2.42%  "

$ws.Range("E30").Value = "  -4.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("E34").Value = "  +0.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.097.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.502"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.90%  "

$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.719.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0510"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.406"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("E51").Value = "  -0.27%  "

# --- Rows 37-39: coin ranking reshuffled (PaxDollar, MXToken and VeChain swapped rank) ---
$ws.Range("B37").Value = "PaxDollar"
$ws.Range("C37").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.27%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0152"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "

